$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate data row for "Carbon dioxide, non-fossil" /
# "air::non-urban air or from high stacks" (row 8), shifting all
# subsequent rows up by one.
$ws.Rows.Item(8).Delete()
